# Generate Report for Handback
# - Marks the two files as handed back (status + handback datetime)
# - Adds "Latest Target File" / "Latest Handback File" hyperlink columns
#   (F, G) to the zh-cn and de-de detail sheets for both rows.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Status text: every cell that used to read "Ready for handoff" now
# reads the handback status (Overview rollup columns + per-language
# Status column on both detail sheets).
# ---------------------------------------------------------------------
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# Latest Handback DateTime (column H): zh-cn got handed back a little
# before de-de.
# ---------------------------------------------------------------------
$zhcn.Range("H2").Value = "2016-03-18 16:53:33"
$zhcn.Range("H3").Value = "2016-03-18 16:53:33"

$dede.Range("H2").Value = "2016-03-18 16:53:38"
$dede.Range("H3").Value = "2016-03-18 16:53:38"

# ---------------------------------------------------------------------
# New columns F (Latest Target File) / G (Latest Handback File):
# same source-doc / translated-file links as columns A and D, now
# also exposed after handback. The whole Hyperlinks collection for
# each sheet is rebuilt (delete + re-add in reading order) so the
# relationship ids come out in the same left-to-right, top-to-bottom
# order Excel itself would produce.
# ---------------------------------------------------------------------

# zh-cn: row 2 = 4e02d25f..., row 3 = fe0e1dab...
$zhcn.Range("F2").Value = "4e02d25f-208e-478f-8b11-99c7ff49e0f5.md"
$zhcn.Range("G2").Value = "4e02d25f-208e-478f-8b11-99c7ff49e0f5.6d900131576a7630fa91e1b65db76e878e57c41c.zh-cn.xlf"
$zhcn.Range("F3").Value = "fe0e1dab-8ca0-4760-bb1e-5f3ccc084aa5.md"
$zhcn.Range("G3").Value = "fe0e1dab-8ca0-4760-bb1e-5f3ccc084aa5.48da20e1fad582794a40babfa9f822ac02be869b.zh-cn.xlf"

$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/cbd5dd9054f9f2ee2c2c636c1c3674bca88a65a4/e2e/4e02d25f-208e-478f-8b11-99c7ff49e0f5.md", "", "", "4e02d25f-208e-478f-8b11-99c7ff49e0f5.md")
$zhcn.Hyperlinks.Add($zhcn.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/cbd5dd9054f9f2ee2c2c636c1c3674bca88a65a4/e2e/4e02d25f-208e-478f-8b11-99c7ff49e0f5.md", "", "", ".md")
$zhcn.Hyperlinks.Add($zhcn.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f401d565e66aaf70fa6f6aa1c4bb77d916bc74a9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4e02d25f-208e-478f-8b11-99c7ff49e0f5.6d900131576a7630fa91e1b65db76e878e57c41c.zh-cn.xlf", "", "", "4e02d25f-208e-478f-8b11-99c7ff49e0f5.6d900131576a7630fa91e1b65db76e878e57c41c.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/cbd5dd9054f9f2ee2c2c636c1c3674bca88a65a4/e2e/4e02d25f-208e-478f-8b11-99c7ff49e0f5.md", "", "", "4e02d25f-208e-478f-8b11-99c7ff49e0f5.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f401d565e66aaf70fa6f6aa1c4bb77d916bc74a9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4e02d25f-208e-478f-8b11-99c7ff49e0f5.6d900131576a7630fa91e1b65db76e878e57c41c.zh-cn.xlf", "", "", "4e02d25f-208e-478f-8b11-99c7ff49e0f5.6d900131576a7630fa91e1b65db76e878e57c41c.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/cbd5dd9054f9f2ee2c2c636c1c3674bca88a65a4/e2e/fe0e1dab-8ca0-4760-bb1e-5f3ccc084aa5.md", "", "", "fe0e1dab-8ca0-4760-bb1e-5f3ccc084aa5.md")
$zhcn.Hyperlinks.Add($zhcn.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/cbd5dd9054f9f2ee2c2c636c1c3674bca88a65a4/e2e/fe0e1dab-8ca0-4760-bb1e-5f3ccc084aa5.md", "", "", ".md")
$zhcn.Hyperlinks.Add($zhcn.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f401d565e66aaf70fa6f6aa1c4bb77d916bc74a9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/fe0e1dab-8ca0-4760-bb1e-5f3ccc084aa5.48da20e1fad582794a40babfa9f822ac02be869b.zh-cn.xlf", "", "", "fe0e1dab-8ca0-4760-bb1e-5f3ccc084aa5.48da20e1fad582794a40babfa9f822ac02be869b.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/cbd5dd9054f9f2ee2c2c636c1c3674bca88a65a4/e2e/fe0e1dab-8ca0-4760-bb1e-5f3ccc084aa5.md", "", "", "fe0e1dab-8ca0-4760-bb1e-5f3ccc084aa5.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f401d565e66aaf70fa6f6aa1c4bb77d916bc74a9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/fe0e1dab-8ca0-4760-bb1e-5f3ccc084aa5.48da20e1fad582794a40babfa9f822ac02be869b.zh-cn.xlf", "", "", "fe0e1dab-8ca0-4760-bb1e-5f3ccc084aa5.48da20e1fad582794a40babfa9f822ac02be869b.zh-cn.xlf")

# de-de: row 2 = 4e02d25f..., row 3 = fe0e1dab...
$dede.Range("F2").Value = "4e02d25f-208e-478f-8b11-99c7ff49e0f5.md"
$dede.Range("G2").Value = "4e02d25f-208e-478f-8b11-99c7ff49e0f5.6d900131576a7630fa91e1b65db76e878e57c41c.de-de.xlf"
$dede.Range("F3").Value = "fe0e1dab-8ca0-4760-bb1e-5f3ccc084aa5.md"
$dede.Range("G3").Value = "fe0e1dab-8ca0-4760-bb1e-5f3ccc084aa5.48da20e1fad582794a40babfa9f822ac02be869b.de-de.xlf"

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/cbd5dd9054f9f2ee2c2c636c1c3674bca88a65a4/e2e/4e02d25f-208e-478f-8b11-99c7ff49e0f5.md", "", "", "4e02d25f-208e-478f-8b11-99c7ff49e0f5.md")
$dede.Hyperlinks.Add($dede.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/cbd5dd9054f9f2ee2c2c636c1c3674bca88a65a4/e2e/4e02d25f-208e-478f-8b11-99c7ff49e0f5.md", "", "", ".md")
$dede.Hyperlinks.Add($dede.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5c804321798627df453182e4fa95a8221eada2a9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4e02d25f-208e-478f-8b11-99c7ff49e0f5.6d900131576a7630fa91e1b65db76e878e57c41c.de-de.xlf", "", "", "4e02d25f-208e-478f-8b11-99c7ff49e0f5.6d900131576a7630fa91e1b65db76e878e57c41c.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/cbd5dd9054f9f2ee2c2c636c1c3674bca88a65a4/e2e/4e02d25f-208e-478f-8b11-99c7ff49e0f5.md", "", "", "4e02d25f-208e-478f-8b11-99c7ff49e0f5.md")
$dede.Hyperlinks.Add($dede.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5c804321798627df453182e4fa95a8221eada2a9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4e02d25f-208e-478f-8b11-99c7ff49e0f5.6d900131576a7630fa91e1b65db76e878e57c41c.de-de.xlf", "", "", "4e02d25f-208e-478f-8b11-99c7ff49e0f5.6d900131576a7630fa91e1b65db76e878e57c41c.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/cbd5dd9054f9f2ee2c2c636c1c3674bca88a65a4/e2e/fe0e1dab-8ca0-4760-bb1e-5f3ccc084aa5.md", "", "", "fe0e1dab-8ca0-4760-bb1e-5f3ccc084aa5.md")
$dede.Hyperlinks.Add($dede.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/cbd5dd9054f9f2ee2c2c636c1c3674bca88a65a4/e2e/fe0e1dab-8ca0-4760-bb1e-5f3ccc084aa5.md", "", "", ".md")
$dede.Hyperlinks.Add($dede.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5c804321798627df453182e4fa95a8221eada2a9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/fe0e1dab-8ca0-4760-bb1e-5f3ccc084aa5.48da20e1fad582794a40babfa9f822ac02be869b.de-de.xlf", "", "", "fe0e1dab-8ca0-4760-bb1e-5f3ccc084aa5.48da20e1fad582794a40babfa9f822ac02be869b.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/cbd5dd9054f9f2ee2c2c636c1c3674bca88a65a4/e2e/fe0e1dab-8ca0-4760-bb1e-5f3ccc084aa5.md", "", "", "fe0e1dab-8ca0-4760-bb1e-5f3ccc084aa5.md")
$dede.Hyperlinks.Add($dede.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5c804321798627df453182e4fa95a8221eada2a9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/fe0e1dab-8ca0-4760-bb1e-5f3ccc084aa5.48da20e1fad582794a40babfa9f822ac02be869b.de-de.xlf", "", "", "fe0e1dab-8ca0-4760-bb1e-5f3ccc084aa5.48da20e1fad582794a40babfa9f822ac02be869b.de-de.xlf")

Write-Host "Handback report generated."
